# Weekly price update: insert a new record as row 6, pushing the
# existing rows 6-8 down to rows 7-9 (dates 2022-11-28, 2021-12-07,
# 2022-03-18 respectively keep their values; only a brand-new weekly
# record for 2022-11-30 is added above them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 6; this shifts the old rows 6,7,8
# down to 7,8,9 and keeps all of their data/formatting intact.
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the new weekly record.
$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(6, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(6, 4).Value = 44895
$ws.Cells.Item(6, 5).Value = 15
$ws.Cells.Item(6, 6).Value = 100112030
$ws.Cells.Item(6, 7).Value = "Poroto granado"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 200
$ws.Cells.Item(6, 11).Value = 1200
$ws.Cells.Item(6, 12).Value = 1300
$ws.Cells.Item(6, 13).Value = 1255
$ws.Cells.Item(6, 14).Value = "$/kilo"
$ws.Cells.Item(6, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(6, 16).Value = 1255
$ws.Cells.Item(6, 17).Value = 1
$ws.Cells.Item(6, 18).Value = "Hortaliza"
